$wb = $excel.ActiveWorkbook
$ws5 = $wb.Worksheets.Item(5)
$ws5.Activate()
$win = $excel.ActiveWindow
Write-Host "ScrollRow:" $win.ScrollRow
Write-Host "ScrollColumn:" $win.ScrollColumn
